# Refresh crypto-symbol snapshot (Price / Volume(1h) columns, plus a
# BOLO / CoinbaseStockToken row swap) to match the latest scrape.
# Leading "'" forces literal text so numeric-/percent-looking values stay
# as plain text (matching the sheet's existing inline-string cells)
# instead of being auto-converted to numbers/percentages by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'331.84"
$ws.Range("E2").Value = "'1.10%"
$ws.Range("D3").Value = "'44.66"
$ws.Range("E3").Value = "'1.48%"
$ws.Range("D4").Value = "'5.549"
$ws.Range("E4").Value = "'-0.52%"
$ws.Range("D5").Value = "'0.08224"
$ws.Range("E5").Value = "'2.44%"
$ws.Range("D6").Value = "'2.036"
$ws.Range("E6").Value = "'2.85%"
$ws.Range("D7").Value = "'0.9769"
$ws.Range("E7").Value = "'3.00%"
$ws.Range("D8").Value = "'0.1124"
$ws.Range("E8").Value = "'-3.58%"
$ws.Range("D9").Value = "'0.1904"
$ws.Range("E9").Value = "'2.71%"
$ws.Range("D10").Value = "'10.25"
$ws.Range("E10").Value = "'-13.36%"
$ws.Range("D11").Value = "'0.09981"
$ws.Range("E11").Value = "'2.45%"
$ws.Range("D12").Value = "'0.04671"
$ws.Range("E12").Value = "'-1.06%"
$ws.Range("E13").Value = "'-0.73%"
$ws.Range("E14").Value = "'-1.68%"
$ws.Range("D15").Value = "'0.04108"
$ws.Range("E15").Value = "'-2.76%"
$ws.Range("D16").Value = "'0.006004"
$ws.Range("E16").Value = "'0.55%"
$ws.Range("D17").Value = "'3.364"
$ws.Range("E17").Value = "'-0.15%"
$ws.Range("D18").Value = "'4.439"
$ws.Range("E18").Value = "'2.34%"
$ws.Range("D19").Value = "'2.612"
$ws.Range("E19").Value = "'2.48%"
$ws.Range("E20").Value = "'-3.55%"
$ws.Range("D21").Value = "'0.1371"
$ws.Range("E21").Value = "'-2.18%"
$ws.Range("D22").Value = "'0.2492"
$ws.Range("E22").Value = "'-0.80%"
$ws.Range("D23").Value = "'0.001303"
$ws.Range("E23").Value = "'3.97%"
$ws.Range("D24").Value = "'0.004412"
$ws.Range("E24").Value = "'2.55%"
$ws.Range("D25").Value = "'0.0001281"
$ws.Range("E25").Value = "'7.34%"
$ws.Range("D26").Value = "'0.0003742"
$ws.Range("E26").Value = "'-5.95%"
$ws.Range("D38").Value = "'0.02779"
$ws.Range("E38").Value = "'7.08%"
$ws.Range("D39").Value = "'0.05727"
$ws.Range("E39").Value = "'3.62%"
$ws.Range("D40").Value = "'0.007643"
$ws.Range("E40").Value = "'0.87%"
$ws.Range("D41").Value = "'0.1422"
$ws.Range("E41").Value = "'1.64%"
$ws.Range("D42").Value = "'0.007541"
$ws.Range("E42").Value = "'-2.14%"
$ws.Range("D43").Value = "'0.001974"
$ws.Range("D44").Value = "'0.008322"
$ws.Range("E44").Value = "'-0.71%"
$ws.Range("E45").Value = "'-1.06%"
$ws.Range("E46").Value = "'-0.19%"
$ws.Range("D47").Value = "'0.0005804"
$ws.Range("E47").Value = "'-0.13%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.002522"
$ws.Range("E48").Value = "'9.54%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.003720"
$ws.Range("E49").Value = "'-23.14%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.19%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.19%"
